# Daily cryptos-list refresh (prices + 1h volume %) generated by GitHub Actions.
# Rows 28/29 and 48/49 also swap rank order (BinanceUSD/EthereumClassic and
# RenderToken/BabyDogeCoin respectively), so Coin name + Link are rewritten too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value. Values that parse as plain numbers are given a leading
# apostrophe (a doubled '' in this single-quoted PowerShell literal = one literal
# quote char) so Excel keeps them as text -- matching the Price column's existing
# text values such as '27.039.53' or '1.00' -- instead of coercing them to a number.
$updates = [ordered]@{
    'D2' = '27.007.29'
    'E2' = '  +2.25%  '
    'D3' = '1.654.99'
    'E3' = '  +2.91%  '
    'E4' = '  +0.05%  '
    'D5' = '''215.23'
    'D6' = '''0.512'
    'E6' = '  +2.27%  '
    'E7' = '  -0.01%  '
    'D8' = '''0.250'
    'E8' = '  +2.28%  '
    'E9' = '  +1.38%  '
    'D10' = '''20.22'
    'E10' = '  +4.68%  '
    'E11' = '  +2.30%  '
    'D12' = '1.890.04'
    'E12' = '  +2.94%  '
    'D13' = '1.647.88'
    'E13' = '  +2.59%  '
    'E14' = '  +1.97%  '
    'E15' = '  +2.20%  '
    'D16' = '''65.20'
    'E16' = '  +2.55%  '
    'D17' = '27.013.42'
    'E17' = '  +2.25%  '
    'D18' = '''236.12'
    'E18' = '  +0.92%  '
    'D19' = '0.0₃0731'
    'E19' = '  +0.71%  '
    'D20' = '''7.72'
    'E20' = '  +0.15%  '
    'E21' = '  +0.03%  '
    'D22' = '''4.42'
    'E22' = '  +3.46%  '
    'D23' = '''9.34'
    'E23' = '  +3.97%  '
    'D24' = '''2.20'
    'E24' = '  +2.65%  '
    'D25' = '''145.55'
    'E25' = '  -1.10%  '
    'D26' = '''7.10'
    'E26' = '  +1.63%  '
    'E27' = '  +0.76%  '
    'B28' = 'BinanceUSD'
    'C28' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'D28' = '''1.00'
    'E28' = '  -0.06%  '
    'B29' = 'EthereumClassic'
    'C29' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D29' = '''15.82'
    'E29' = '  +2.23%  '
    'E30' = '  +0.17%  '
    'E31' = '  +1.31%  '
    'D32' = '1.553.52'
    'E32' = '  +4.18%  '
    'D33' = '''3.32'
    'E33' = '  +2.90%  '
    'E34' = '  +4.73%  '
    'D35' = '''1.60'
    'E35' = '  +8.18%  '
    'E36' = '  -0.07%  '
    'E37' = '  +3.35%  '
    'D38' = '''0.894'
    'E38' = '  +8.52%  '
    'E39' = '  +2.76%  '
    'E40' = '  +3.19%  '
    'D42' = '''2.25'
    'E42' = '  +2.67%  '
    'D43' = '''65.57'
    'E43' = '  +7.57%  '
    'D44' = '1.796.74'
    'E44' = '  +2.77%  '
    'E45' = '  +1.81%  '
    'D46' = '''0.916'
    'E46' = '  -2.40%  '
    'D47' = '''90.24'
    'E47' = '  +1.08%  '
    'B48' = 'RenderToken'
    'C48' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D48' = '''1.52'
    'E48' = '  +1.47%  '
    'B49' = 'BabyDogeCoin'
    'C49' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D49' = '0.0₆0104'
    'E49' = '  +12.09%  '
    'D50' = '''0.0985'
    'E50' = '  +2.27%  '
    'E51' = '  +0.86%  '
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
